$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain their text formatting so numeric-looking
# strings (e.g. "1.30", "7.00", "0.320") are not coerced into numbers and
# lose significant trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '70.507.80'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '3.559.65'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '612.57'
$ws.Range('E5').Value = '  +5.64%  '
$ws.Range('D6').Value = '172.84'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('E7').Value = '  +1.52%  '
$ws.Range('D8').Value = '3.554.88'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.196'
$ws.Range('E10').Value = '  +3.75%  '
$ws.Range('D11').Value = '7.48'
$ws.Range('E11').Value = '  +11.32%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '46.64'
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '4.140.47'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '8.38'
$ws.Range('E16').Value = '  -2.03%  '
$ws.Range('D17').Value = '616.28'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '3.562.84'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').Value = '70.664.92'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').Value = '0.882'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').Value = '9.39'
$ws.Range('E23').Value = '  -16.45%  '
$ws.Range('D24').Value = '15.92'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '97.16'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '2.62'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('D29').Value = '33.45'
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('D30').Value = '9.03'
$ws.Range('E30').Value = '  -3.16%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').Value = '3.05'
$ws.Range('E32').Value = '  -3.42%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').Value = '1.30'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '7.00'
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').Value = '575.32'
$ws.Range('E35').Value = '  -9.06%  '
$ws.Range('D36').Value = '3.70'
$ws.Range('E36').Value = '  +5.48%  '
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('D38').Value = '10.81'
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0476'
$ws.Range('E39').Value = '  +4.28%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '57.31'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +3.96%  '
$ws.Range('D43').Value = '3.378.79'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '0.320'
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('D45').Value = '33.19'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').Value = '2.98'
$ws.Range('E46').Value = '  +7.66%  '
$ws.Range('D47').Value = '0.0₃0706'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('D48').Value = '2.62'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('D50').Value = '133.83'
$ws.Range('E50').Value = '  +1.73%  '
